$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is the same (45189) for
# every data row (rows 2-171). The edit bumps that value by one day (45190)
# across the whole column range.
$range = $ws.Range("C2:C171")
$range.Value = 45190
